# regen sval data to filter save games
# Update the B:E (and derived G=sum) columns for rows 2-16 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B,C,D,E,G values per row (F "Win" column is left untouched).
$data = @{
    2  = @(1.445647641019636,   1.626987699542094,   0.1496068669990043,  0.5333859586016987, 3.755628166162433)
    3  = @(0.2881169905109251,  0.04103571897497393, 0.1496068669990043,  0.5333859586016987, 1.012145535086602)
    4  = @(3.272327238179451,   1.626987699542094,   0.1496068669990043,  0.5333859586016987, 5.582307763322248)
    5  = @(3.272327238179451,   1.626987699542094,   0.1496068669990043,  0.5333859586016987, 5.582307763322248)
    6  = @(0.2881169905109251,  0.3048912486333797,  189.6080260415259,   0.5333859586016987, 190.7344202392719)
    7  = @(3.272327238179451,   1.626987699542094,   0.7210945179870265,  0.5333859586016987, 6.15379541431027)
    8  = @(0.04172184405617529, 0.04103571897497393, 0.7210945179870265,  0.5333859586016987, 1.337238039619874)
    9  = @(3.272327238179451,   1.626987699542094,   0.1496068669990043,  0.5333859586016987, 5.582307763322248)
    10 = @(0.1169995834814548,  0.3048912486333797,  0.1496068669990043,  0.5333859586016987, 1.104883657715537)
    11 = @(3.272327238179451,   1.626987699542094,   0.7210945179870265,  0.5333859586016987, 6.15379541431027)
    12 = @(0.01253208636536152, 1.626987699542094,   3.223369029078222,   0.5333859586016987, 5.396274773587376)
    13 = @(1.445647641019636,   1.626987699542094,   3.223369029078222,   0.5333859586016987, 6.82939032824165)
    14 = @(3.272327238179451,   1.626987699542094,   3.223369029078222,   0.5333859586016987, 8.656069925401464)
    15 = @(3.272327238179451,   1.626987699542094,   0.7210945179870265,  0.5333859586016987, 6.15379541431027)
    16 = @(3.272327238179451,   0.3048912486333797,  0.1496068669990043,  0.5333859586016987, 4.260211312413533)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
